$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# New column M (2021) mirrors the formatting of column L (2020) cell-for-cell.
# Header row 2 (thin bottom-border blank cell) gets extended into column M.
$ws.Range("L2").Copy() | Out-Null
$ws.Range("M2").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("L3").Copy() | Out-Null
$ws.Range("M3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("M3").Value = 2021

$ws.Range("L4").Copy() | Out-Null
$ws.Range("M4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("M4").Value = 952

$ws.Range("L5").Copy() | Out-Null
$ws.Range("M5").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("M5").Value = 10437

$ws.Range("L6").Copy() | Out-Null
$ws.Range("M6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("M6").Value = 2253

$ws.Range("L7").Copy() | Out-Null
$ws.Range("M7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("M7").Value = 8184

$ws.Range("L8").Copy() | Out-Null
$ws.Range("M8").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("M8").Value = 14020

$ws.Range("L9").Copy() | Out-Null
$ws.Range("M9").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("M9").Value = 5139

$ws.Range("L10").Copy() | Out-Null
$ws.Range("M10").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("M10").Value = 8881

$excel.CutCopyMode = $false

# Update the saved selection to match the authored state
$ws.Range("P8").Select()
